$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain plain text (matches the
# original inline-string cell type); force text format first so Excel does not
# auto-convert them to numbers, then restore the cell's original (Normal) style.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D15", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value2 = '43.219.49'
$ws.Range("E2").Value2 = '  +1.63%  '
$ws.Range("D3").Value2 = '2.367.49'
$ws.Range("E3").Value2 = '  +6.46%  '
$ws.Range("E4").Value2 = '  -0.18%  '
$ws.Range("D5").Value2 = '308.32'
$ws.Range("E5").Value2 = '  +3.29%  '
$ws.Range("D6").Value2 = '105.84'
$ws.Range("E6").Value2 = '  -6.02%  '
$ws.Range("E7").Value2 = '  +1.39%  '
$ws.Range("E8").Value2 = '  -0.25%  '
$ws.Range("D9").Value2 = '0.637'
$ws.Range("E9").Value2 = '  +3.98%  '
$ws.Range("D10").Value2 = '42.73'
$ws.Range("E10").Value2 = '  -5.75%  '
$ws.Range("D11").Value2 = '0.0940'
$ws.Range("E11").Value2 = '  +1.20%  '
$ws.Range("D12").Value2 = '8.95'
$ws.Range("E12").Value2 = '  +0.45%  '
$ws.Range("E13").Value2 = '  +11.10%  '
$ws.Range("E14").Value2 = '  +1.20%  '
$ws.Range("D15").Value2 = '16.54'
$ws.Range("E15").Value2 = '  +8.91%  '
$ws.Range("D16").Value2 = '2.724.87'
$ws.Range("E16").Value2 = '  +6.45%  '
$ws.Range("D17").Value2 = '2.368.58'
$ws.Range("E17").Value2 = '  +5.90%  '
$ws.Range("D18").Value2 = '43.139.27'
$ws.Range("E18").Value2 = '  +1.75%  '
$ws.Range("D19").Value2 = '7.42'
$ws.Range("E19").Value2 = '  +0.80%  '
$ws.Range("E20").Value2 = '  +1.99%  '
$ws.Range("D21").Value2 = '75.47'
$ws.Range("E21").Value2 = '  +2.24%  '
$ws.Range("D22").Value2 = '3.41'
$ws.Range("E22").Value2 = '  -3.54%  '
$ws.Range("E23").Value2 = '  +8.70%  '
$ws.Range("D24").Value2 = '252.32'
$ws.Range("E24").Value2 = '  +9.74%  '
$ws.Range("D25").Value2 = '8.93'
$ws.Range("E25").Value2 = '  -5.36%  '
$ws.Range("D26").Value2 = '11.98'
$ws.Range("E26").Value2 = '  +1.66%  '
$ws.Range("E27").Value2 = '  +0.06%  '
$ws.Range("B28").Value2 = 'Toncoin'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value2 = '2.25'
$ws.Range("E28").Value2 = '  +1.20%  '
$ws.Range("B29").Value2 = 'InjectiveProtocol'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value2 = '38.79'
$ws.Range("E29").Value2 = '  -0.63%  '
$ws.Range("D30").Value2 = '22.63'
$ws.Range("E30").Value2 = '  +6.83%  '
$ws.Range("D31").Value2 = '173.15'
$ws.Range("E31").Value2 = '  -0.69%  '
$ws.Range("E32").Value2 = '  -2.18%  '
$ws.Range("E33").Value2 = '  +2.15%  '
$ws.Range("D34").Value2 = '5.87'
$ws.Range("E34").Value2 = '  +2.22%  '
$ws.Range("D35").Value2 = '0.132'
$ws.Range("E35").Value2 = '  +3.44%  '
$ws.Range("D36").Value2 = '4.95'
$ws.Range("E36").Value2 = '  +0.18%  '
$ws.Range("D37").Value2 = '0.0377'
$ws.Range("E37").Value2 = '  +2.01%  '
$ws.Range("D38").Value2 = '4.02'
$ws.Range("E38").Value2 = '  -6.78%  '
$ws.Range("E39").Value2 = '  +0.17%  '
$ws.Range("D40").Value2 = '2.78'
$ws.Range("E40").Value2 = '  +11.20%  '
$ws.Range("D41").Value2 = '1.52'
$ws.Range("E41").Value2 = '  +15.05%  '
$ws.Range("D42").Value2 = '72.00'
$ws.Range("E42").Value2 = '  +1.32%  '
$ws.Range("D43").Value2 = '0.231'
$ws.Range("E43").Value2 = '  -3.00%  '
$ws.Range("E44").Value2 = '  -0.16%  '
$ws.Range("D45").Value2 = '12.24'
$ws.Range("E45").Value2 = '  -7.05%  '
$ws.Range("D46").Value2 = '5.66'
$ws.Range("E46").Value2 = '  +2.59%  '
$ws.Range("D47").Value2 = '9.37'
$ws.Range("E47").Value2 = '  +9.80%  '
$ws.Range("D48").Value2 = '112.56'
$ws.Range("E48").Value2 = '  +7.13%  '
$ws.Range("E49").Value2 = '  -3.41%  '
$ws.Range("E50").Value2 = '  +0.78%  '
$ws.Range("D51").Value2 = '1.495.61'
$ws.Range("E51").Value2 = '  +4.22%  '

foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
